$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G50: append two new narrower-vocab ids to existing broader list ---
$ws.Range('G50').Value = 'vocab:1000,vocab:1038,vocab:1034,vocab:1028,vocab:1064,vocab:1057,vocab:1252,vocab:1244,vocab:1263,vocab:1259,vocab:1278,vocab:1172,vocab:1291,vocab:1299'

# --- B118:B192 label text changes (drop trailing " matrix") ---
$ws.Range('B118').Value = 'air'
$ws.Range('B119').Value = 'sediment'
$ws.Range('B120').Value = 'soil'
$ws.Range('B121').Value = 'water'
$ws.Range('B122').Value = 'biota'
$ws.Range('B123').Value = 'deposition'
$ws.Range('B124').Value = 'whole air'
$ws.Range('B125').Value = 'gas phase'
$ws.Range('B126').Value = 'TSP'
$ws.Range('B127').Value = 'PM 10'
$ws.Range('B128').Value = 'PM 2.5'
$ws.Range('B129').Value = 'PM 5'
$ws.Range('B130').Value = 'surface'
$ws.Range('B131').Value = 'sewage sludge'
$ws.Range('B132').Value = 'stratified'
$ws.Range('B133').Value = 'surface'
$ws.Range('B134').Value = 'topsoil'
$ws.Range('B135').Value = 'subsoil A'
$ws.Range('B136').Value = 'subsoil B'
$ws.Range('B137').Value = 'horizon A'
$ws.Range('B138').Value = 'horizon B'
$ws.Range('B139').Value = 'horizon C'
$ws.Range('B140').Value = 'horizon E'
$ws.Range('B141').Value = 'horizon G'
$ws.Range('B142').Value = 'horizon M'
$ws.Range('B143').Value = 'horizon O'
$ws.Range('B144').Value = 'coastal'
$ws.Range('B145').Value = 'brackish'
$ws.Range('B146').Value = 'ocean'
$ws.Range('B147').Value = 'surface freshwater'
$ws.Range('B148').Value = 'surface freshwater - reservoir'
$ws.Range('B149').Value = 'surface freshwater - river'
$ws.Range('B150').Value = 'ground water'
$ws.Range('B151').Value = 'fish'
$ws.Range('B152').Value = 'barley'
$ws.Range('B153').Value = 'grass'
$ws.Range('B154').Value = 'hop'
$ws.Range('B155').Value = 'moss'
$ws.Range('B156').Value = 'mustard'
$ws.Range('B157').Value = 'maize'
$ws.Range('B158').Value = 'oat'
$ws.Range('B159').Value = 'peat moss'
$ws.Range('B160').Value = 'pine'
$ws.Range('B161').Value = 'potato'
$ws.Range('B162').Value = 'rape'
$ws.Range('B163').Value = 'soya'
$ws.Range('B164').Value = 'sunflover'
$ws.Range('B165').Value = 'spruce'
$ws.Range('B166').Value = 'triticale'
$ws.Range('B167').Value = 'wheat'
$ws.Range('B168').Value = 'dry'
$ws.Range('B169').Value = 'mixed'
$ws.Range('B170').Value = 'wet'
$ws.Range('B171').Value = 'blood-whole blood '
$ws.Range('B172').Value = 'blood -plasma '
$ws.Range('B173').Value = 'blood -serum '
$ws.Range('B174').Value = 'cord blood-whole blood '
$ws.Range('B175').Value = 'cord blood-plasma '
$ws.Range('B176').Value = 'cord blood-serum '
$ws.Range('B177').Value = 'urine-spot '
$ws.Range('B178').Value = 'urine-24h '
$ws.Range('B179').Value = 'urine-morning urine '
$ws.Range('B180').Value = 'saliva and/or sputum '
$ws.Range('B181').Value = 'semen '
$ws.Range('B182').Value = 'hair'
$ws.Range('B183').Value = 'exhaled breath condensate'
$ws.Range('B184').Value = 'red blood cells'
$ws.Range('B185').Value = 'breast milk'
$ws.Range('B186').Value = 'adipose tissue/fat'
$ws.Range('B187').Value = 'all toe nails'
$ws.Range('B188').Value = 'big toe nails'
$ws.Range('B189').Value = 'dermal wipes'
$ws.Range('B190').Value = 'amniotic fluid'
$ws.Range('B191').Value = 'placenta tissue'
$ws.Range('B192').Value = 'human'

# --- Y column: set dct:modified date (forced to Text so date stays literal string) ---
$dateCells = @('Y118','Y119','Y171','Y172','Y173','Y174','Y175','Y176','Y177','Y178','Y179','Y180','Y181','Y182')
foreach ($c in $dateCells) {
    $ws.Range($c).NumberFormat = "@"
    $ws.Range($c).Value = "2024-01-10"
}
# --- New rows 308-323: additional vocabulary terms added in this revision ---
$ws.Range('A308').Value = 'vocab:1288'
$ws.Range('B308').Value = 'exhaled air'
$ws.Range('D308').Value = 'EA'
$ws.Range('G308').Value = 'vocab:1172'

$ws.Range('A309').Value = 'vocab:1289'
$ws.Range('B309').Value = 'dried blood spots'
$ws.Range('D309').Value = 'DBS'
$ws.Range('G309').Value = 'vocab:1172'

$ws.Range('A310').Value = 'vocab:1290'
$ws.Range('B310').Value = 'volumetric absorptive microsample'
$ws.Range('D310').Value = 'VAMS'
$ws.Range('G310').Value = 'vocab:1172'

$ws.Range('A311').Value = 'vocab:1291'
$ws.Range('B311').Value = 'environmental matrix'
$ws.Range('D311').Value = 'SWB'

$ws.Range('A312').Value = 'vocab:1292'
$ws.Range('B312').Value = 'silicone wrist band'
$ws.Range('D312').Value = 'IDUST'
$ws.Range('G312').Value = 'vocab:1291'

$ws.Range('A313').Value = 'vocab:1293'
$ws.Range('B313').Value = 'indoor dust'
$ws.Range('D313').Value = 'ODUST'
$ws.Range('G313').Value = 'vocab:1291'

$ws.Range('A314').Value = 'vocab:1294'
$ws.Range('B314').Value = 'outdoor dust'
$ws.Range('D314').Value = 'IAIR'
$ws.Range('G314').Value = 'vocab:1291'

$ws.Range('A315').Value = 'vocab:1295'
$ws.Range('B315').Value = 'indoor air-stationary'
$ws.Range('D315').Value = 'OAIR'
$ws.Range('G315').Value = 'vocab:1291'

$ws.Range('A316').Value = 'vocab:1296'
$ws.Range('B316').Value = 'outdoor air-stationary'
$ws.Range('D316').Value = 'PAIR'
$ws.Range('G316').Value = 'vocab:1291'

$ws.Range('A317').Value = 'vocab:1297'
$ws.Range('B317').Value = 'personal air'
$ws.Range('D317').Value = 'SW'
$ws.Range('G317').Value = 'vocab:1291'

$ws.Range('A318').Value = 'vocab:1298'
$ws.Range('B318').Value = 'surface wipe'
$ws.Range('G318').Value = 'vocab:1291'

$ws.Range('A319').Value = 'vocab:1299'
$ws.Range('B319').Value = 'collection place'

$ws.Range('A320').Value = 'vocab:1300'
$ws.Range('B320').Value = 'participants home'
$ws.Range('G320').Value = 'vocab:1299'

$ws.Range('A321').Value = 'vocab:1301'
$ws.Range('B321').Value = 'educational establishment'
$ws.Range('G321').Value = 'vocab:1299'

$ws.Range('A322').Value = 'vocab:1302'
$ws.Range('B322').Value = 'workplace'
$ws.Range('G322').Value = 'vocab:1299'

$ws.Range('A323').Value = 'vocab:1303'
$ws.Range('B323').Value = 'healthcare establishment'
$ws.Range('G323').Value = 'vocab:1299'

# --- Row 324: fully blank data row, only the modified-date stamp is set ---

# --- Y308:Y324 modified-date stamps for all newly appended rows ---
$newDateCells = @('Y308','Y309','Y310','Y311','Y312','Y313','Y314','Y315','Y316','Y317','Y318','Y319','Y320','Y321','Y322','Y323','Y324')
foreach ($c in $newDateCells) {
    $ws.Range($c).NumberFormat = "@"
    $ws.Range($c).Value = "2024-01-10"
}
